$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.861952666666666
$ws.Range("H2").Value = 14.585858
$ws.Range("I2").Value = 0.3995648519435639
$ws.Range("J2").Value = 0.3995648519435638
$ws.Range("M2").Value = 45.90594266666667
$ws.Range("N2").Value = 137.717828
$ws.Range("O2").Value = 0.3954672001633582
$ws.Range("P2").Value = 0.3954672001633583
$ws.Range("Q2").Value = 223.1925203640471
$ws.Range("R2").Value = 2008.732683276424
$ws.Range("S2").Value = 0.158014793281808
$ws.Range("T2").Value = 0.158014793281808

$ws.Range("G3").Value = 4.861952666666666
$ws.Range("H3").Value = 14.585858
$ws.Range("I3").Value = 0.3995648519435639
$ws.Range("J3").Value = 0.3995648519435638
$ws.Range("M3").Value = 40.44578266666667
$ws.Range("N3").Value = 121.337348
$ws.Range("O3").Value = 0.3484294080560655
$ws.Range("P3").Value = 0.3484294080560656
$ws.Range("Q3").Value = 196.6454808916204
$ws.Range("R3").Value = 1769.809328024584
$ws.Range("S3").Value = 0.1392201448427054
$ws.Range("T3").Value = 0.1392201448427055

$ws.Range("G4").Value = 4.861952666666666
$ws.Range("H4").Value = 14.585858
$ws.Range("I4").Value = 0.3995648519435639
$ws.Range("J4").Value = 0.3995648519435638
$ws.Range("M4").Value = 12.761795
$ws.Range("N4").Value = 38.28538500000001
$ws.Range("O4").Value = 0.1099393900775594
$ws.Range("P4").Value = 0.1099393900775594
$ws.Range("Q4").Value = 62.04724323170333
$ws.Range("R4").Value = 558.42518908533
$ws.Range("S4").Value = 0.04392791611910574
$ws.Range("T4").Value = 0.04392791611910574

$ws.Range("G5").Value = 4.861952666666666
$ws.Range("H5").Value = 14.585858
$ws.Range("I5").Value = 0.3995648519435639
$ws.Range("J5").Value = 0.3995648519435638
$ws.Range("M5").Value = 16.966758
$ws.Range("N5").Value = 50.900274
$ws.Range("O5").Value = 0.1461640017030168
$ws.Range("P5").Value = 0.1461640017030168
$ws.Range("Q5").Value = 82.491574302788
$ws.Range("R5").Value = 742.4241687250918
$ws.Range("S5").Value = 0.05840199769994473
$ws.Range("T5").Value = 0.05840199769994473

$ws.Range("G6").Value = 0.6001993333333334
$ws.Range("H6").Value = 1.800598
$ws.Range("I6").Value = 0.04932556406896855
$ws.Range("J6").Value = 0.04932556406896854
$ws.Range("M6").Value = 45.90594266666667
$ws.Range("N6").Value = 137.717828
$ws.Range("O6").Value = 0.3954672001633582
$ws.Range("P6").Value = 0.3954672001633583
$ws.Range("Q6").Value = 27.55271618457156
$ws.Range("R6").Value = 247.974445661144
$ws.Range("S6").Value = 0.01950664271883333
$ws.Range("T6").Value = 0.01950664271883333

$ws.Range("G7").Value = 0.6001993333333334
$ws.Range("H7").Value = 1.800598
$ws.Range("I7").Value = 0.04932556406896855
$ws.Range("J7").Value = 0.04932556406896854
$ws.Range("M7").Value = 40.44578266666667
$ws.Range("N7").Value = 121.337348
$ws.Range("O7").Value = 0.3484294080560655
$ws.Range("P7").Value = 0.3484294080560656
$ws.Range("Q7").Value = 24.27553179267823
$ws.Range("R7").Value = 218.479786134104
$ws.Range("S7").Value = 0.01718647709058225
$ws.Range("T7").Value = 0.01718647709058225

$ws.Range("G8").Value = 0.6001993333333334
$ws.Range("H8").Value = 1.800598
$ws.Range("I8").Value = 0.04932556406896855
$ws.Range("J8").Value = 0.04932556406896854
$ws.Range("M8").Value = 12.761795
$ws.Range("N8").Value = 38.28538500000001
$ws.Range("O8").Value = 0.1099393900775594
$ws.Range("P8").Value = 0.1099393900775594
$ws.Range("Q8").Value = 7.659620851136668
$ws.Range("R8").Value = 68.93658766023002
$ws.Range("S8").Value = 0.005422822428973981
$ws.Range("T8").Value = 0.005422822428973982

$ws.Range("G9").Value = 0.6001993333333334
$ws.Range("H9").Value = 1.800598
$ws.Range("I9").Value = 0.04932556406896855
$ws.Range("J9").Value = 0.04932556406896854
$ws.Range("M9").Value = 16.966758
$ws.Range("N9").Value = 50.900274
$ws.Range("O9").Value = 0.1461640017030168
$ws.Range("P9").Value = 0.1461640017030168
$ws.Range("Q9").Value = 10.183436840428
$ws.Range("R9").Value = 91.650931563852
$ws.Range("S9").Value = 0.007209621830578983
$ws.Range("T9").Value = 0.007209621830578983

$ws.Range("G10").Value = 4.206754333333333
$ws.Range("H10").Value = 12.620263
$ws.Range("I10").Value = 0.3457193616641432
$ws.Range("J10").Value = 0.3457193616641432
$ws.Range("M10").Value = 45.90594266666667
$ws.Range("N10").Value = 137.717828
$ws.Range("O10").Value = 0.3954672001633582
$ws.Range("P10").Value = 0.3954672001633583
$ws.Range("Q10").Value = 193.1150232387515
$ws.Range("R10").Value = 1738.035209148764
$ws.Range("S10").Value = 0.1367206679995822
$ws.Range("T10").Value = 0.1367206679995822

$ws.Range("G11").Value = 4.206754333333333
$ws.Range("H11").Value = 12.620263
$ws.Range("I11").Value = 0.3457193616641432
$ws.Range("J11").Value = 0.3457193616641432
$ws.Range("M11").Value = 40.44578266666667
$ws.Range("N11").Value = 121.337348
$ws.Range("O11").Value = 0.3484294080560655
$ws.Range("P11").Value = 0.3484294080560656
$ws.Range("Q11").Value = 170.1454714980582
$ws.Range("R11").Value = 1531.309243482524
$ws.Range("S11").Value = 0.1204587925381583
$ws.Range("T11").Value = 0.1204587925381583

$ws.Range("G12").Value = 4.206754333333333
$ws.Range("H12").Value = 12.620263
$ws.Range("I12").Value = 0.3457193616641432
$ws.Range("J12").Value = 0.3457193616641432
$ws.Range("M12").Value = 12.761795
$ws.Range("N12").Value = 38.28538500000001
$ws.Range("O12").Value = 0.1099393900775594
$ws.Range("P12").Value = 0.1099393900775594
$ws.Range("Q12").Value = 53.68573641736166
$ws.Range("R12").Value = 483.171627756255
$ws.Range("S12").Value = 0.03800817575935908
$ws.Range("T12").Value = 0.03800817575935909

$ws.Range("G13").Value = 4.206754333333333
$ws.Range("H13").Value = 12.620263
$ws.Range("I13").Value = 0.3457193616641432
$ws.Range("J13").Value = 0.3457193616641432
$ws.Range("M13").Value = 16.966758
$ws.Range("N13").Value = 50.900274
$ws.Range("O13").Value = 0.1461640017030168
$ws.Range("P13").Value = 0.1461640017030168
$ws.Range("Q13").Value = 71.37498273911798
$ws.Range("R13").Value = 642.3748446520618
$ws.Range("S13").Value = 0.05053172536704371
$ws.Range("T13").Value = 0.05053172536704372

$ws.Range("G14").Value = 2.499212666666667
$ws.Range("H14").Value = 7.497638
$ws.Range("I14").Value = 0.2053902223233243
$ws.Range("J14").Value = 0.2053902223233243
$ws.Range("M14").Value = 45.90594266666667
$ws.Range("N14").Value = 137.717828
$ws.Range("O14").Value = 0.3954672001633582
$ws.Range("P14").Value = 0.3954672001633583
$ws.Range("Q14").Value = 114.7287133878071
$ws.Range("R14").Value = 1032.558420490264
$ws.Range("S14").Value = 0.08122509616313474
$ws.Range("T14").Value = 0.08122509616313477

$ws.Range("G15").Value = 2.499212666666667
$ws.Range("H15").Value = 7.497638
$ws.Range("I15").Value = 0.2053902223233243
$ws.Range("J15").Value = 0.2053902223233243
$ws.Range("M15").Value = 40.44578266666667
$ws.Range("N15").Value = 121.337348
$ws.Range("O15").Value = 0.3484294080560655
$ws.Range("P15").Value = 0.3484294080560656
$ws.Range("Q15").Value = 101.0826123537804
$ws.Range("R15").Value = 909.743511184024
$ws.Range("S15").Value = 0.0715639935846196
$ws.Range("T15").Value = 0.07156399358461961

$ws.Range("G16").Value = 2.499212666666667
$ws.Range("H16").Value = 7.497638
$ws.Range("I16").Value = 0.2053902223233243
$ws.Range("J16").Value = 0.2053902223233243
$ws.Range("M16").Value = 12.761795
$ws.Range("N16").Value = 38.28538500000001
$ws.Range("O16").Value = 0.1099393900775594
$ws.Range("P16").Value = 0.1099393900775594
$ws.Range("Q16").Value = 31.89443971340334
$ws.Range("R16").Value = 287.0499574206301
$ws.Range("S16").Value = 0.0225804757701206
$ws.Range("T16").Value = 0.02258047577012061

$ws.Range("G17").Value = 2.499212666666667
$ws.Range("H17").Value = 7.497638
$ws.Range("I17").Value = 0.2053902223233243
$ws.Range("J17").Value = 0.2053902223233243
$ws.Range("M17").Value = 16.966758
$ws.Range("N17").Value = 50.900274
$ws.Range("O17").Value = 0.1461640017030168
$ws.Range("P17").Value = 0.1461640017030168
$ws.Range("Q17").Value = 42.403536505868
$ws.Range("R17").Value = 381.631828552812
$ws.Range("S17").Value = 0.03002065680544938
$ws.Range("T17").Value = 0.03002065680544938
